$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-10-30T07:01:42.353576+00:00"
$ws.Range("K3").Value = "2025-10-30T07:01:44.631832+00:00"
$ws.Range("K4").Value = "2025-10-30T07:01:44.631894+00:00"
$ws.Range("K5").Value = "2025-10-30T07:01:44.631922+00:00"
$ws.Range("K6").Value = "2025-10-30T07:01:44.631943+00:00"
$ws.Range("K7").Value = "2025-10-30T07:01:44.631962+00:00"
$ws.Range("K8").Value = "2025-10-30T07:01:47.365292+00:00"
$ws.Range("K9").Value = "2025-10-30T07:01:47.365324+00:00"
$ws.Range("K10").Value = "2025-10-30T07:01:50.218952+00:00"
$ws.Range("K11").Value = "2025-10-30T07:01:52.535673+00:00"
$ws.Range("K12").Value = "2025-10-30T07:01:52.535706+00:00"
$ws.Range("K13").Value = "2025-10-30T07:01:54.998708+00:00"
$ws.Range("K14").Value = "2025-10-30T07:01:54.998740+00:00"
$ws.Range("K15").Value = "2025-10-30T07:01:54.998758+00:00"
$ws.Range("K16").Value = "2025-10-30T07:01:54.998776+00:00"
$ws.Range("K17").Value = "2025-10-30T07:02:02.330791+00:00"
$ws.Range("K18").Value = "2025-10-30T07:02:05.183984+00:00"
$ws.Range("K19").Value = "2025-10-30T07:02:07.866335+00:00"
$ws.Range("K20").Value = "2025-10-30T07:02:10.535743+00:00"
$ws.Range("K21").Value = "2025-10-30T07:02:10.535775+00:00"
$ws.Range("K22").Value = "2025-10-30T07:02:10.535793+00:00"
$ws.Range("K23").Value = "2025-10-30T07:02:12.875695+00:00"
$ws.Range("K24").Value = "2025-10-30T07:02:12.875725+00:00"
$ws.Range("K25").Value = "2025-10-30T07:02:12.875743+00:00"
$ws.Range("K26").Value = "2025-10-30T07:02:12.875759+00:00"
$ws.Range("K27").Value = "2025-10-30T07:02:12.875776+00:00"
$ws.Range("K28").Value = "2025-10-30T07:02:20.579502+00:00"
$ws.Range("K29").Value = "2025-10-30T07:02:20.579533+00:00"
$ws.Range("K30").Value = "2025-10-30T07:02:20.579552+00:00"
$ws.Range("K31").Value = "2025-10-30T07:02:20.579570+00:00"
$ws.Range("K32").Value = "2025-10-30T07:02:23.251964+00:00"
$ws.Range("K33").Value = "2025-10-30T07:02:23.251995+00:00"
$ws.Range("K34").Value = "2025-10-30T07:02:23.252014+00:00"
$ws.Range("K35").Value = "2025-10-30T07:02:25.478586+00:00"
$ws.Range("K36").Value = "2025-10-30T07:02:25.478617+00:00"
$ws.Range("K37").Value = "2025-10-30T07:02:25.478635+00:00"
$ws.Range("K38").Value = "2025-10-30T07:02:25.478652+00:00"
$ws.Range("K39").Value = "2025-10-30T07:02:25.478669+00:00"
$ws.Range("K40").Value = "2025-10-30T07:02:25.478685+00:00"
$ws.Range("K41").Value = "2025-10-30T07:02:25.478701+00:00"
$ws.Range("K42").Value = "2025-10-30T07:02:25.478723+00:00"
$ws.Range("K43").Value = "2025-10-30T07:02:25.478739+00:00"
$ws.Range("K44").Value = "2025-10-30T07:02:27.669444+00:00"
$ws.Range("K45").Value = "2025-10-30T07:02:27.669476+00:00"
$ws.Range("K46").Value = "2025-10-30T07:02:32.479432+00:00"
$ws.Range("K47").Value = "2025-10-30T07:02:34.715968+00:00"
$ws.Range("K48").Value = "2025-10-30T07:02:34.716007+00:00"
$ws.Range("K49").Value = "2025-10-30T07:02:34.716029+00:00"
$ws.Range("K50").Value = "2025-10-30T07:02:34.716048+00:00"
